$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 255:256, shifting the existing rows 255-276 down to 257-278.
$ws.Range("A255:A256").EntireRow.Insert()

# Populate the newly inserted row 255 (Primera quality, week of 44610).
$ws.Range("A255").Value = 1
$ws.Range("B255").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C255").Value = "Arica y Parinacota"
$ws.Range("D255").Value = 44610
$ws.Range("E255").Value = 15
$ws.Range("F255").Value = 100112032
$ws.Range("G255").Value = "Zapallo italiano"
$ws.Range("H255").Value = "Huracán"
$ws.Range("I255").Value = "Primera"
$ws.Range("J255").Value = 130
$ws.Range("K255").Value = 4500
$ws.Range("L255").Value = 5000
$ws.Range("M255").Value = 4750
$ws.Range("N255").Value = "$/caja 70 unidades"
$ws.Range("O255").Value = "Región de Arica y Parinacota"
$ws.Range("P255").Value = 68
$ws.Range("Q255").Value = 70
$ws.Range("R255").Value = "Hortaliza"

# Populate the newly inserted row 256 (Segunda quality, week of 44610).
$ws.Range("A256").Value = 1
$ws.Range("B256").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C256").Value = "Arica y Parinacota"
$ws.Range("D256").Value = 44610
$ws.Range("E256").Value = 15
$ws.Range("F256").Value = 100112032
$ws.Range("G256").Value = "Zapallo italiano"
$ws.Range("H256").Value = "Huracán"
$ws.Range("I256").Value = "Segunda"
$ws.Range("J256").Value = 160
$ws.Range("K256").Value = 4000
$ws.Range("L256").Value = 4500
$ws.Range("M256").Value = 4250
$ws.Range("N256").Value = "$/caja 100 unidades"
$ws.Range("O256").Value = "Región de Arica y Parinacota"
$ws.Range("P256").Value = 42
$ws.Range("Q256").Value = 100
$ws.Range("R256").Value = "Hortaliza"
